# [MOSIP-14369] Fix: boolean values
#
# The upload template's "is_active" column (D2:D5) held the values as the
# formula =TRUE() (a numeric 1). That's wrong for a seed/reference data
# file that's consumed as plain text - it should just contain the literal
# text "TRUE". Replace the formulas with literal text "TRUE" values,
# keeping the existing text format/style, then update the selection to
# reflect the edited range (D2:D5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..5) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.ClearContents()
    $cell.NumberFormat = "@"
    # Write as a formula first so Excel doesn't auto-coerce the literal
    # "TRUE" text into a Boolean, then convert it to a plain value in
    # place (paste-special values) so the stored cell is literal text,
    # not a formula - matching the target content.
    $cell.Formula = '="TRUE"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

$ws.Range("D2:D5").Select()
